$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 29.82164133333333
$ws.Range("H2").Value = 89.464924
$ws.Range("I2").Value = 0.02335016309719764
$ws.Range("J2").Value = 0.02335016309719765
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.699506666666667
$ws.Range("N2").Value = 14.09852
$ws.Range("O2").Value = 0.9660495246229048
$ws.Range("P2").Value = 0.9660495246229047
$ws.Range("Q2").Value = 140.1470022569422
$ws.Range("R2").Value = 1261.32302031248
$ws.Range("S2").Value = 0.02255741395991507
$ws.Range("T2").Value = 0.02255741395991508

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 29.82164133333333
$ws.Range("H3").Value = 89.464924
$ws.Range("I3").Value = 0.02335016309719764
$ws.Range("J3").Value = 0.02335016309719765
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.03395047537709522
$ws.Range("P3").Value = 0.03395047537709522
$ws.Range("Q3").Value = 4.925272698783556
$ws.Range("R3").Value = 44.327454289052
$ws.Range("S3").Value = 0.0007927491372825659
$ws.Range("T3").Value = 0.0007927491372825661

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1160.126729666667
$ws.Range("H4").Value = 3480.380189
$ws.Range("I4").Value = 0.9083721465342723
$ws.Range("J4").Value = 0.9083721465342726
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.699506666666667
$ws.Range("N4").Value = 14.09852
$ws.Range("O4").Value = 0.9660495246229048
$ws.Range("P4").Value = 0.9660495246229047
$ws.Range("Q4").Value = 5452.023300246698
$ws.Range("R4").Value = 49068.20970222027
$ws.Range("S4").Value = 0.8775324803401214
$ws.Range("T4").Value = 0.8775324803401215

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1160.126729666667
$ws.Range("H5").Value = 3480.380189
$ws.Range("I5").Value = 0.9083721465342723
$ws.Range("J5").Value = 0.9083721465342726
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.03395047537709522
$ws.Range("P5").Value = 0.03395047537709522
$ws.Range("Q5").Value = 191.6038237093774
$ws.Range("R5").Value = 1724.434413384397
$ws.Range("S5").Value = 0.03083966619415094
$ws.Range("T5").Value = 0.03083966619415095

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2607933333333334
$ws.Range("H6").Value = 0.7823800000000001
$ws.Range("I6").Value = 0.0002041995878070102
$ws.Range("J6").Value = 0.0002041995878070102
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.699506666666667
$ws.Range("N6").Value = 14.09852
$ws.Range("O6").Value = 0.9660495246229048
$ws.Range("P6").Value = 0.9660495246229047
$ws.Range("Q6").Value = 1.225600008622223
$ws.Range("R6").Value = 11.0304000776
$ws.Range("S6").Value = 0.0001972669147291553
$ws.Range("T6").Value = 0.0001972669147291553

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2607933333333334
$ws.Range("H7").Value = 0.7823800000000001
$ws.Range("I7").Value = 0.0002041995878070102
$ws.Range("J7").Value = 0.0002041995878070102
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.03395047537709522
$ws.Range("P7").Value = 0.03395047537709522
$ws.Range("Q7").Value = 0.04307201841555557
$ws.Range("R7").Value = 0.3876481657400001
$ws.Range("S7").Value = 0.000006932673077854892
$ws.Range("T7").Value = 0.000006932673077854893

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Serping1"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 86.94000199999999
$ws.Range("H8").Value = 260.820006
$ws.Range("I8").Value = 0.0680734907807228
$ws.Range("J8").Value = 0.06807349078072282
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.699506666666667
$ws.Range("N8").Value = 14.09852
$ws.Range("O8").Value = 0.9660495246229048
$ws.Range("P8").Value = 0.9660495246229047
$ws.Range("Q8").Value = 408.5751189990133
$ws.Range("R8").Value = 3677.17607099112
$ws.Range("S8").Value = 0.06576236340813896
$ws.Range("T8").Value = 0.06576236340813897

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Serping1"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 86.94000199999999
$ws.Range("H9").Value = 260.820006
$ws.Range("I9").Value = 0.0680734907807228
$ws.Range("J9").Value = 0.06807349078072282
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.03395047537709522
$ws.Range("P9").Value = 0.03395047537709522
$ws.Range("Q9").Value = 14.35880787031533
$ws.Range("R9").Value = 129.229270832838
$ws.Range("S9").Value = 0.002311127372583848
$ws.Range("T9").Value = 0.002311127372583848
